# Refresh the "Created Room ID" values (column L) in the booking dataset
# sheet with a new batch of booking IDs produced by a fresh test run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (cell, new "Created Room ID" value) pairs, row by row.
$updates = @(
    @{ Cell = "L2"; Value = "501726448" }
    @{ Cell = "L9"; Value = "501739339" }
    @{ Cell = "L10"; Value = "501742110" }
    @{ Cell = "L11"; Value = "501744782" }
    @{ Cell = "L12"; Value = "501747262" }
    @{ Cell = "L13"; Value = "501750036" }
    @{ Cell = "L14"; Value = "501752825" }
    @{ Cell = "L15"; Value = "501755118" }
    @{ Cell = "L16"; Value = "501757539" }
    @{ Cell = "L17"; Value = "501760306" }
    @{ Cell = "L18"; Value = "501763047" }
    @{ Cell = "L19"; Value = "501765110" }
    @{ Cell = "L20"; Value = "501767473" }
    @{ Cell = "L21"; Value = "501769640" }
    @{ Cell = "L22"; Value = "501772146" }
    @{ Cell = "L23"; Value = "501774837" }
    @{ Cell = "L24"; Value = "501778314" }
    @{ Cell = "L25"; Value = "501781234" }
    @{ Cell = "L26"; Value = "501784796" }
    @{ Cell = "L27"; Value = "501789032" }
    @{ Cell = "L28"; Value = "501792672" }
    @{ Cell = "L29"; Value = "501796195" }
    @{ Cell = "L30"; Value = "501799602" }
    @{ Cell = "L31"; Value = "501803087" }
    @{ Cell = "L32"; Value = "501806235" }
    @{ Cell = "L33"; Value = "501809237" }
    @{ Cell = "L34"; Value = "501810977" }
    @{ Cell = "L35"; Value = "501813223" }
    @{ Cell = "L36"; Value = "501815477" }
    @{ Cell = "L37"; Value = "501818155" }
    @{ Cell = "L38"; Value = "501820380" }
    @{ Cell = "L39"; Value = "501822604" }
    @{ Cell = "L40"; Value = "501824833" }
    @{ Cell = "L41"; Value = "501827115" }
    @{ Cell = "L42"; Value = "501829269" }
    @{ Cell = "L43"; Value = "501831473" }
    @{ Cell = "L44"; Value = "501833725" }
    @{ Cell = "L45"; Value = "501836195" }
    @{ Cell = "L46"; Value = "501838393" }
    @{ Cell = "L47"; Value = "501840979" }
    @{ Cell = "L48"; Value = "501845769" }
    @{ Cell = "L49"; Value = "501849918" }
    @{ Cell = "L50"; Value = "501859139" }
    @{ Cell = "L51"; Value = "501863171" }
    @{ Cell = "L52"; Value = "501867580" }
    @{ Cell = "L53"; Value = "501854646" }
    @{ Cell = "L54"; Value = "501872414" }
    @{ Cell = "L55"; Value = "501876607" }
    @{ Cell = "L56"; Value = "501881075" }
    @{ Cell = "L57"; Value = "501885683" }
    @{ Cell = "L58"; Value = "501892058" }
    @{ Cell = "L59"; Value = "501896543" }
    @{ Cell = "L60"; Value = "501899839" }
    @{ Cell = "L61"; Value = "501955025" }
    @{ Cell = "L62"; Value = "501960768" }
    @{ Cell = "L63"; Value = "501965491" }
    @{ Cell = "L64"; Value = "501970378" }
    @{ Cell = "L65"; Value = "501974916" }
    @{ Cell = "L66"; Value = "501979576" }
    @{ Cell = "L67"; Value = "501984114" }
    @{ Cell = "L68"; Value = "501989016" }
    @{ Cell = "L69"; Value = "501902720" }
    @{ Cell = "L70"; Value = "501920581" }
    @{ Cell = "L71"; Value = "501927085" }
    @{ Cell = "L72"; Value = "501906988" }
    @{ Cell = "L74"; Value = "501914292" }
    @{ Cell = "L75"; Value = "501934756" }
    @{ Cell = "L76"; Value = "501941139" }
    @{ Cell = "L77"; Value = "501948292" }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    # Force text storage so the numeric-looking ID round-trips as a string
    # (matching the existing shared-string / t="s" cells) instead of becoming
    # a numeric cell, then restore the original "General" number format so the
    # cell style is left untouched.
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.NumberFormat = "General"
}
